# Updates cell values across multiple sheets per the scheduled-runner market-price refresh.
# Each row keeps its Leve Item ID (column G) as the join key; only columns H-N (price/profit
# columns) are refreshed here, matching the upstream data pull.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 674.5
$ws.Range("J18").Value = 648.5
$ws.Range("L18").Value = 648.5
$ws.Range("N18").Value = -1216.5
# Row 70
$ws.Range("H70").Value = 2748.875
$ws.Range("J70").Value = 3012.7144
$ws.Range("L70").Value = 9038.143199999999
$ws.Range("N70").Value = -9578.143199999999
# Row 73
$ws.Range("H73").Value = 2748.875
$ws.Range("J73").Value = 3012.7144
$ws.Range("L73").Value = 9038.143199999999
$ws.Range("N73").Value = -10910.1432
# Row 82
$ws.Range("H82").Value = 4542.222
$ws.Range("I82").Value = 2701.8572
$ws.Range("K82").Value = 8105.571599999999
$ws.Range("M82").Value = -7699.571599999999
# Row 85
$ws.Range("H85").Value = 4542.222
$ws.Range("I85").Value = 2701.8572
$ws.Range("K85").Value = 8105.571599999999
$ws.Range("M85").Value = -6701.571599999999
# Row 100
$ws.Range("H100").Value = 3304
$ws.Range("I100").Value = 2508.8
$ws.Range("J100").Value = 4099.2
$ws.Range("K100").Value = 2508.8
$ws.Range("L100").Value = 4099.2
$ws.Range("M100").Value = -1967.8
$ws.Range("N100").Value = -5181.2
# Row 106
$ws.Range("H106").Value = 6833.6924
$ws.Range("I106").Value = 3530.6365
$ws.Range("K106").Value = 3530.6365
$ws.Range("M106").Value = -2899.6365
# Row 135
$ws.Range("H135").Value = 1314.65
$ws.Range("I135").Value = 1288
$ws.Range("K135").Value = 11592
$ws.Range("M135").Value = -9057
# Row 138
$ws.Range("H138").Value = 3778.9546
$ws.Range("I138").Value = 3787.476
$ws.Range("K138").Value = 11362.428
$ws.Range("M138").Value = -6222.428

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 4160.4
$ws.Range("I2").Value = 4229.5713
$ws.Range("K2").Value = 4229.5713
$ws.Range("M2").Value = -4116.5713
# Row 32
$ws.Range("H32").Value = 27190.023
$ws.Range("I32").Value = 27679
$ws.Range("K32").Value = 27679
$ws.Range("M32").Value = -27392
# Row 45
$ws.Range("H45").Value = 3206.9473
$ws.Range("I45").Value = 1399.125
$ws.Range("J45").Value = 4521.727
$ws.Range("K45").Value = 1399.125
$ws.Range("L45").Value = 4521.727
$ws.Range("M45").Value = -1022.125
$ws.Range("N45").Value = -5275.727
# Row 64
$ws.Range("H64").Value = 249999
$ws.Range("J64").Value = 249999
$ws.Range("L64").Value = 249999
$ws.Range("N64").Value = -250495
# Row 67
$ws.Range("H67").Value = 249999
$ws.Range("J67").Value = 249999
$ws.Range("L67").Value = 249999
$ws.Range("N67").Value = -251715
# Row 116
$ws.Range("H116").Value = 4160.4
$ws.Range("I116").Value = 4229.5713
$ws.Range("K116").Value = 4229.5713
$ws.Range("M116").Value = -1935.5713

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 4160.4
$ws.Range("I3").Value = 4229.5713
$ws.Range("K3").Value = 4229.5713
$ws.Range("M3").Value = -4115.5713
# Row 20
$ws.Range("H20").Value = 5267298.5
$ws.Range("I20").Value = 9093853
$ws.Range("K20").Value = 9093853
$ws.Range("M20").Value = -9093606
# Row 53
$ws.Range("H53").Value = 74999
$ws.Range("J53").Value = 74999
$ws.Range("L53").Value = 74999
$ws.Range("N53").Value = -76147
# Row 127
$ws.Range("H127").Value = 69999
$ws.Range("J127").Value = 69999
$ws.Range("L127").Value = 69999
$ws.Range("N127").Value = -79919

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 7678.7856
$ws.Range("J31").Value = 9956.143
$ws.Range("L31").Value = 9956.143
$ws.Range("N31").Value = -10546.143
# Row 34
$ws.Range("H34").Value = 7678.7856
$ws.Range("J34").Value = 9956.143
$ws.Range("L34").Value = 9956.143
$ws.Range("N34").Value = -10360.143
# Row 44
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = $null
$ws.Range("N44").Value = $null
# Row 62
$ws.Range("H62").Value = 4386.625
$ws.Range("I62").Value = 3799.3333
$ws.Range("J62").Value = 4739
$ws.Range("K62").Value = 3799.3333
$ws.Range("L62").Value = 4739
$ws.Range("M62").Value = -3175.3333
$ws.Range("N62").Value = -5987
# Row 65
$ws.Range("H65").Value = 4386.625
$ws.Range("I65").Value = 3799.3333
$ws.Range("J65").Value = 4739
$ws.Range("K65").Value = 18996.6665
$ws.Range("L65").Value = 23695
$ws.Range("M65").Value = -15876.6665
$ws.Range("N65").Value = -29935

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 378
$ws.Range("I5").Value = 378
$ws.Range("K5").Value = 1134
$ws.Range("M5").Value = -1022
# Row 135
$ws.Range("H135").Value = 378
$ws.Range("I135").Value = 378
$ws.Range("K135").Value = 3402
$ws.Range("M135").Value = -867

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 2909.8262
$ws.Range("I122").Value = 2260.2354
$ws.Range("J122").Value = 4750.3335
$ws.Range("K122").Value = 6780.706200000001
$ws.Range("L122").Value = 14251.0005
$ws.Range("M122").Value = -4330.706200000001
$ws.Range("N122").Value = -19151.0005
# Row 134
$ws.Range("H134").Value = 35000
$ws.Range("J134").Value = 35000
$ws.Range("L134").Value = 105000
$ws.Range("N134").Value = -110070
# Row 136
$ws.Range("H136").Value = 25222.5
$ws.Range("J136").Value = 25222.5
$ws.Range("L136").Value = 75667.5
$ws.Range("N136").Value = -80767.5
# Row 138
$ws.Range("H138").Value = 84999
$ws.Range("J138").Value = 84999
$ws.Range("L138").Value = 84999
$ws.Range("N138").Value = -95279
# Row 139
$ws.Range("H139").Value = 26773.334
$ws.Range("J139").Value = 26773.334
$ws.Range("L139").Value = 26773.334
$ws.Range("N139").Value = -37053.334

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 48
$ws.Range("H48").Value = 28178
$ws.Range("I48").Value = 22267.25
$ws.Range("K48").Value = 22267.25
$ws.Range("M48").Value = -21606.25
# Row 100
$ws.Range("H100").Value = 2992.5293
$ws.Range("I100").Value = 2427
$ws.Range("K100").Value = 2427
$ws.Range("M100").Value = -1886
# Row 122
$ws.Range("H122").Value = 4446.0835
$ws.Range("I122").Value = 3195
$ws.Range("J122").Value = 5339.7144
$ws.Range("K122").Value = 9585
$ws.Range("L122").Value = 16019.1432
$ws.Range("M122").Value = -7135
$ws.Range("N122").Value = -20919.1432
# Row 132
$ws.Range("H132").Value = 74399.06
$ws.Range("I132").Value = 83079.60000000001
$ws.Range("K132").Value = 249238.8
$ws.Range("M132").Value = -246708.8
# Row 134
$ws.Range("H134").Value = 90665
$ws.Range("J134").Value = 84500
$ws.Range("L134").Value = 84500
$ws.Range("N134").Value = -94640
# Row 137
$ws.Range("H137").Value = 88000
$ws.Range("J137").Value = 88000
$ws.Range("L137").Value = 88000
$ws.Range("N137").Value = -98200

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 58
$ws.Range("H58").Value = 18750
$ws.Range("I58").Value = 18333.334
$ws.Range("J58").Value = 20000
$ws.Range("K58").Value = 18333.334
$ws.Range("L58").Value = 20000
$ws.Range("M58").Value = -18025.334
$ws.Range("N58").Value = -20616
# Row 76
$ws.Range("H76").Value = 24996
$ws.Range("I76").Value = 24996
$ws.Range("K76").Value = 24996
$ws.Range("M76").Value = -24681
# Row 79
$ws.Range("H79").Value = 24996
$ws.Range("I79").Value = 24996
$ws.Range("K79").Value = 24996
$ws.Range("M79").Value = -23904
# Row 81
$ws.Range("H81").Value = 1831.6666
$ws.Range("I81").Value = 1091.1111
$ws.Range("K81").Value = 2182.2222
$ws.Range("M81").Value = -1121.2222
# Row 84
$ws.Range("H84").Value = 1831.6666
$ws.Range("I84").Value = 1091.1111
$ws.Range("K84").Value = 10911.111
$ws.Range("M84").Value = -5607.111000000001
# Row 136
$ws.Range("H136").Value = 2062.8
$ws.Range("I136").Value = 1724
$ws.Range("K136").Value = 5172
$ws.Range("M136").Value = -2622
